# Applies the "write excel summary sheet" commit:
#  - adds 9 new (empty) worksheets after the existing "汇总信息" sheet
#  - appends rows 10-21 of new data (with merges) to sheet 1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheets, in order, after all existing sheets.
# ---------------------------------------------------------------------------
$newSheetNames = @(
    "含钢量汇总",
    "计算参数",
    "周期",
    "内力",
    "位移角",
    "整体验算结果",
    "楼层分布数据",
    "调整系数",
    "工程量"
)

foreach ($name in $newSheetNames) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
    $newSheet.Name = $name
}

# ---------------------------------------------------------------------------
# 2. Populate the new rows (10-21) on the first sheet ("汇总信息").
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item(1)

# row 10
$ws.Range("A10").Value = "质量"
$ws.Range("C10").Value = "活载质量"
$ws.Range("D10").Value = 23980.664

# row 11
$ws.Range("C11").Value = "恒载质量"
$ws.Range("D11").Value = 198539.547
$ws.Range("F11").Value = 222520.219

# row 12
$ws.Range("A12").Value = "层间位移角"
$ws.Range("B12").Value = "风荷载"
$ws.Range("C12").Value = "X向"
$ws.Range("D12").Value = 878
$ws.Range("E12").Value = "楼层"
$ws.Range("F12").Value = 53

# row 13
$ws.Range("C13").Value = "Y向"
$ws.Range("D13").Value = 469
$ws.Range("E13").Value = "楼层"
$ws.Range("F13").Value = 53

# row 14
$ws.Range("B14").Value = "-偏心"
$ws.Range("C14").Value = "X向"
$ws.Range("D14").Value = 646
$ws.Range("E14").Value = "楼层"
$ws.Range("F14").Value = 53

# row 15
$ws.Range("C15").Value = "Y向"
$ws.Range("D15").Value = 1.37
$ws.Range("E15").Value = "楼层"
$ws.Range("F15").Value = 53

# row 16
$ws.Range("B16").Value = "限值"
$ws.Range("D16").Value = 500

# row 17
$ws.Range("A17").Value = "位移比"
$ws.Range("B17").Value = "+偏心"
$ws.Range("C17").Value = "X向"
$ws.Range("D17").Value = 1.44
$ws.Range("E17").Value = "楼层"
$ws.Range("F17").Value = 53

# row 18
$ws.Range("C18").Value = "Y向"
$ws.Range("D18").Value = 1.37
$ws.Range("E18").Value = "楼层"
$ws.Range("F18").Value = 53

# row 19
$ws.Range("C19").Value = "X向"
$ws.Range("D19").Value = 1.44
$ws.Range("E19").Value = "楼层"
$ws.Range("F19").Value = 53

# row 20
$ws.Range("C20").Value = "Y向"
$ws.Range("E20").Value = "楼层"

# row 21
$ws.Range("B21").Value = "限值"
$ws.Range("D21").Value = "1.2 or 1.4"

# ---------------------------------------------------------------------------
# 3. Merge the cells that belong together.
# ---------------------------------------------------------------------------
$ws.Range("A10:B11").Merge()
$ws.Range("A12:A16").Merge()
$ws.Range("B12:B13").Merge()
$ws.Range("B14:B15").Merge()
$ws.Range("B16:C16").Merge()
$ws.Range("D16:F16").Merge()
$ws.Range("A17:A21").Merge()
$ws.Range("B17:B18").Merge()
$ws.Range("B19:B20").Merge()
$ws.Range("B21:C21").Merge()
$ws.Range("D21:F21").Merge()

$ws.Activate()
